# màj couleur saison 22-23
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the U15M match info on row 17 (previously blank)
$ws.Range("B17").Value = "JEUDI"
$ws.Range("C17").Value = "15H"
$ws.Range("D17").Value = "TOTO"
$ws.Range("E17").Value = "MONTAIGNE"

# Move the active selection to match the saved workbook state
$ws.Range("E18").Select()
